$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.396.61"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.265.91"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.70%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "614.88"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.78"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.68%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.267.05"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.76%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.162"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.81"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.497"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.69%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "39.04"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.76%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.794.83"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.462.61"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.45"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.259.93"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.40%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "505.95"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.80%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.757"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +3.68%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.13"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.70"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "87.09"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.05"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.22%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.72%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.40"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.128"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +45.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.02"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.88"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "28.02"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.98%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.48"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.34%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +18.51%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0790"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +15.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "494.85"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0423"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.73%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.85"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.77%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.294"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.006.08"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +6.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "29.00"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.32%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +5.89%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.38%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.25%  "
